# Apply edits described in commit "adding data up to 15th"
# - Adds Puerto Rico (AR) figures for the week of 29 Jun-4 Jul (rows 142-147)
# - Corrects a handful of already-published figures for 05 07 2020 / 06 07 2020 (rows 156-157)
# - Fills in the data rows for 06 07 2020 and 07 07 2020 (rows 158-159, dates already present)
# - Appends new daily rows for 08 07 2020 .. 15 07 2020 (rows 160-167)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Back-fill Puerto Rico (column AR) for rows 142-147 ---
$ws.Cells.Item(142, 44).Value = 9.6638655
$ws.Cells.Item(143, 44).Value = 14.1026846
$ws.Cells.Item(144, 44).Value = 12.5000302
$ws.Cells.Item(145, 44).Value = 6.3433735
$ws.Cells.Item(146, 44).Value = 9.9206349
$ws.Cells.Item(147, 44).Value = 12.3188934

# --- Small value corrections on rows 156 and 157 ---
$ws.Cells.Item(156, 6).Value = 32.4653096
$ws.Cells.Item(157, 6).Value = 32.7789865
$ws.Cells.Item(157, 20).Value = 20.1683719
$ws.Cells.Item(157, 22).Value = 29.4708224
$ws.Cells.Item(157, 24).Value = 12.8346568
$ws.Cells.Item(157, 26).Value = 12.7986214
$ws.Cells.Item(157, 34).Value = 21.213696
$ws.Cells.Item(157, 43).Value = 12.4952696
$ws.Cells.Item(157, 49).Value = 33.705177
$ws.Cells.Item(157, 54).Value = 14.9071715

# --- Fill existing rows 158-159 and append new rows 160-167 (A:BE) ---

# Row 158
$row158 = New-Object 'object[,]' 1,57
$row158[0,0] = "06 07 2020"
$row158[0,1] = 16.9183253
$row158[0,2] = 29.5143513
$row158[0,3] = 24.107177
$row158[0,5] = 33.3704249
$row158[0,6] = 17.6924464
$row158[0,7] = 14.8147658
$row158[0,8] = 10.2229111
$row158[0,9] = 13.9324487
$row158[0,10] = 14.7277228
$row158[0,11] = 26.3150971
$row158[0,12] = 25.6116608
$row158[0,14] = 10.7705287
$row158[0,15] = 24.528655
$row158[0,16] = 25.3590334
$row158[0,17] = 13.4592071
$row158[0,18] = 19.9827343
$row158[0,19] = 20.4922932
$row158[0,20] = 17.6406222
$row158[0,21] = 30.1696393
$row158[0,22] = 10.475409
$row158[0,23] = 13.0648484
$row158[0,24] = 9.9
$row158[0,25] = 13.3700591
$row158[0,26] = 16.4760546
$row158[0,27] = 19.1893983
$row158[0,29] = 28.7638491
$row158[0,30] = 18.4329957
$row158[0,31] = 19.6781178
$row158[0,32] = 19.1976577
$row158[0,33] = 21.2905758
$row158[0,34] = 10.0958565
$row158[0,35] = 11.642863
$row158[0,36] = 16.4989034
$row158[0,37] = 20.0819124
$row158[0,38] = 10.9822438
$row158[0,39] = 15.9457223
$row158[0,40] = 25.3034962
$row158[0,41] = 12.9228472
$row158[0,42] = 12.7566773
$row158[0,44] = 10.8671279
$row158[0,45] = 30.1666944
$row158[0,46] = 21.6263594
$row158[0,47] = 21.9745574
$row158[0,48] = 34.1965975
$row158[0,49] = 24.8026304
$row158[0,50] = 15.4307353
$row158[0,52] = 9.700578500000001
$row158[0,53] = 15.4091442
$row158[0,54] = 17.1123663
$row158[0,55] = 17.1643016
$row158[0,56] = 15.8328005
$ws.Range("A158:BE158").Value = $row158

# Row 159
$row159 = New-Object 'object[,]' 1,57
$row159[0,0] = "07 07 2020"
$row159[0,1] = 18.540096
$row159[0,2] = 30.3507816
$row159[0,3] = 23.9970369
$row159[0,5] = 33.802515
$row159[0,6] = 18.0073822
$row159[0,7] = 14.8983356
$row159[0,8] = 10.1358924
$row159[0,9] = 14.1420118
$row159[0,10] = 15.097561
$row159[0,11] = 26.7837643
$row159[0,12] = 26.2870662
$row159[0,14] = 11.3264192
$row159[0,15] = 24.9419881
$row159[0,16] = 26.7417414
$row159[0,17] = 13.8691768
$row159[0,18] = 20.102531
$row159[0,19] = 21.5423498
$row159[0,20] = 18.6622986
$row159[0,21] = 31.2915365
$row159[0,22] = 10.4951572
$row159[0,23] = 13.3537736
$row159[0,24] = 10.1571114
$row159[0,25] = 13.5461185
$row159[0,26] = 16.5215635
$row159[0,27] = 19.7837845
$row159[0,29] = 29.4474594
$row159[0,30] = 18.7756844
$row159[0,31] = 19.9430736
$row159[0,32] = 18.157946
$row159[0,33] = 21.1176784
$row159[0,34] = 10.3566122
$row159[0,35] = 11.6876809
$row159[0,36] = 17.1188141
$row159[0,37] = 20.2281426
$row159[0,38] = 11.0039296
$row159[0,39] = 16.387481
$row159[0,40] = 25.547398
$row159[0,41] = 13.2311248
$row159[0,42] = 12.6931724
$row159[0,44] = 10.8970687
$row159[0,45] = 30.105139
$row159[0,46] = 21.9915205
$row159[0,47] = 22.8504615
$row159[0,48] = 34.7858815
$row159[0,49] = 25.0153675
$row159[0,50] = 15.6163363
$row159[0,52] = 8.601823
$row159[0,53] = 15.7853958
$row159[0,54] = 17.3608334
$row159[0,55] = 17.542068
$row159[0,56] = 16.6247329
$ws.Range("A159:BE159").Value = $row159

# Row 160
$row160 = New-Object 'object[,]' 1,57
$row160[0,0] = "08 07 2020"
$row160[0,1] = 18.197634
$row160[0,2] = 31.4844674
$row160[0,3] = 24.0797015
$row160[0,5] = 33.7000066
$row160[0,6] = 18.5504864
$row160[0,7] = 15.3646456
$row160[0,8] = 9.551384499999999
$row160[0,9] = 14.0423032
$row160[0,10] = 15.4932735
$row160[0,11] = 27.3962847
$row160[0,12] = 27.3502762
$row160[0,14] = 10.3878116
$row160[0,15] = 25.839906
$row160[0,16] = 27.6209115
$row160[0,17] = 14.2405344
$row160[0,18] = 20.382524
$row160[0,19] = 21.9070784
$row160[0,20] = 18.5532281
$row160[0,21] = 31.7126262
$row160[0,22] = 10.5243153
$row160[0,23] = 13.2871455
$row160[0,24] = 9.926656599999999
$row160[0,25] = 13.7146305
$row160[0,26] = 17.0930184
$row160[0,27] = 19.7885177
$row160[0,29] = 30.3346846
$row160[0,30] = 19.7161342
$row160[0,31] = 19.992823
$row160[0,32] = 18.6846111
$row160[0,33] = 21.4337795
$row160[0,34] = 10.2089805
$row160[0,35] = 11.7189474
$row160[0,36] = 17.1653453
$row160[0,37] = 21.2486593
$row160[0,38] = 11.1026743
$row160[0,39] = 16.6487976
$row160[0,40] = 26.3150892
$row160[0,41] = 13.8597935
$row160[0,42] = 12.9421702
$row160[0,44] = 11.1516148
$row160[0,45] = 30.6113227
$row160[0,46] = 21.1885529
$row160[0,47] = 23.4921564
$row160[0,48] = 34.9544485
$row160[0,49] = 25.3308189
$row160[0,50] = 16.1464001
$row160[0,52] = 9.5851764
$row160[0,53] = 16.2511838
$row160[0,54] = 17.4881269
$row160[0,55] = 17.8434843
$row160[0,56] = 16.8749846
$ws.Range("A160:BE160").Value = $row160

# Row 161
$row161 = New-Object 'object[,]' 1,57
$row161[0,0] = "09 07 2020"
$row161[0,1] = 19.7042641
$row161[0,2] = 32.5005568
$row161[0,3] = 24.4717329
$row161[0,5] = 33.663415
$row161[0,6] = 18.890091
$row161[0,7] = 15.5951172
$row161[0,8] = 9.735512999999999
$row161[0,9] = 14.5314353
$row161[0,10] = 15.7330652
$row161[0,11] = 27.6561341
$row161[0,12] = 28.3518815
$row161[0,14] = 10.5304212
$row161[0,15] = 26.249919
$row161[0,16] = 27.8090574
$row161[0,17] = 14.8367973
$row161[0,18] = 21.1088812
$row161[0,19] = 22.8405003
$row161[0,20] = 18.8234783
$row161[0,21] = 32.50586
$row161[0,22] = 10.8127532
$row161[0,23] = 13.6640493
$row161[0,24] = 10.6006583
$row161[0,25] = 14.2797998
$row161[0,26] = 17.2841197
$row161[0,27] = 20.654161
$row161[0,29] = 30.5350005
$row161[0,30] = 19.9262169
$row161[0,31] = 20.2679215
$row161[0,32] = 18.5492036
$row161[0,33] = 21.4314004
$row161[0,34] = 10.5608046
$row161[0,35] = 11.7918518
$row161[0,36] = 16.8835545
$row161[0,37] = 21.3683338
$row161[0,38] = 11.3579482
$row161[0,39] = 16.8949676
$row161[0,40] = 26.9110247
$row161[0,41] = 14.3032515
$row161[0,42] = 13.2085646
$row161[0,44] = 11.7751507
$row161[0,45] = 30.9650017
$row161[0,46] = 22.6290361
$row161[0,47] = 24.1840717
$row161[0,48] = 35.6002134
$row161[0,49] = 25.122235
$row161[0,50] = 16.1649252
$row161[0,52] = 8.5948767
$row161[0,53] = 16.4536366
$row161[0,54] = 18.0352443
$row161[0,55] = 18.0270042
$row161[0,56] = 17.4189924
$ws.Range("A161:BE161").Value = $row161

# Row 162
$row162 = New-Object 'object[,]' 1,57
$row162[0,0] = "10 07 2020"
$row162[0,1] = 20.382372
$row162[0,2] = 32.4618765
$row162[0,3] = 24.6116366
$row162[0,5] = 34.0362763
$row162[0,6] = 19.1068955
$row162[0,7] = 15.9813237
$row162[0,8] = 9.470913100000001
$row162[0,9] = 14.0736342
$row162[0,10] = 15.7998252
$row162[0,11] = 28.0292812
$row162[0,12] = 29.2316447
$row162[0,14] = 10.6965174
$row162[0,15] = 26.3172584
$row162[0,16] = 28.055521
$row162[0,17] = 15.4020479
$row162[0,18] = 21.6407044
$row162[0,19] = 23.6480803
$row162[0,20] = 19.2725137
$row162[0,21] = 33.4794689
$row162[0,22] = 10.99812
$row162[0,23] = 13.6256903
$row162[0,24] = 10.2540355
$row162[0,25] = 14.5724986
$row162[0,26] = 17.7304814
$row162[0,27] = 20.9043963
$row162[0,29] = 31.3898123
$row162[0,30] = 18.7238945
$row162[0,31] = 20.2710843
$row162[0,32] = 20.8673608
$row162[0,33] = 21.6024612
$row162[0,34] = 10.4267425
$row162[0,35] = 11.7741789
$row162[0,36] = 17.3505791
$row162[0,37] = 22.5339997
$row162[0,38] = 11.5262973
$row162[0,39] = 17.1613543
$row162[0,40] = 27.7904621
$row162[0,41] = 14.2139927
$row162[0,42] = 13.6003307
$row162[0,44] = 11.6738569
$row162[0,45] = 31.6753206
$row162[0,46] = 22.6671
$row162[0,47] = 24.630052
$row162[0,48] = 36.0287863
$row162[0,49] = 25.8646409
$row162[0,50] = 16.479992
$row162[0,52] = 8.480681499999999
$row162[0,53] = 16.9239729
$row162[0,54] = 18.8778053
$row162[0,55] = 18.4857864
$row162[0,56] = 19.0715124
$ws.Range("A162:BE162").Value = $row162

# Row 163
$row163 = New-Object 'object[,]' 1,57
$row163[0,0] = "11 07 2020"
$row163[0,1] = 20.1459391
$row163[0,2] = 33.1856147
$row163[0,3] = 24.9718846
$row163[0,5] = 33.675882
$row163[0,6] = 19.5523072
$row163[0,7] = 16.1401716
$row163[0,8] = 9.980370499999999
$row163[0,9] = 14.0642939
$row163[0,10] = 16.0867698
$row163[0,11] = 28.4399577
$row163[0,12] = 29.8730243
$row163[0,14] = 11.0353535
$row163[0,15] = 26.8472813
$row163[0,16] = 28.4823756
$row163[0,17] = 15.9794238
$row163[0,18] = 21.5582374
$row163[0,19] = 23.9433336
$row163[0,20] = 19.3815299
$row163[0,21] = 33.7468987
$row163[0,22] = 10.8430414
$row163[0,23] = 13.7605188
$row163[0,24] = 10.3998961
$row163[0,25] = 14.7176725
$row163[0,26] = 18.0537978
$row163[0,27] = 21.2060415
$row163[0,29] = 31.7221948
$row163[0,30] = 19.8526649
$row163[0,31] = 20.4488728
$row163[0,32] = 19.7551759
$row163[0,33] = 21.710834
$row163[0,34] = 10.1017042
$row163[0,35] = 11.8543623
$row163[0,36] = 17.7799838
$row163[0,37] = 23.3497166
$row163[0,38] = 11.7677035
$row163[0,39] = 17.6238576
$row163[0,40] = 27.6707333
$row163[0,41] = 14.5636303
$row163[0,42] = 13.675631
$row163[0,44] = 11.8656667
$row163[0,45] = 31.8462533
$row163[0,46] = 23.1889159
$row163[0,47] = 24.8443858
$row163[0,48] = 36.6818059
$row163[0,49] = 25.6927497
$row163[0,50] = 16.6324797
$row163[0,52] = 9.6381801
$row163[0,53] = 17.0988719
$row163[0,54] = 19.4696877
$row163[0,55] = 18.5877909
$row163[0,56] = 19.2370293
$ws.Range("A163:BE163").Value = $row163

# Row 164
$row164 = New-Object 'object[,]' 1,57
$row164[0,0] = "12 07 2020"
$row164[0,1] = 21.0140306
$row164[0,2] = 33.6131484
$row164[0,3] = 25.1359934
$row164[0,5] = 34.0000635
$row164[0,6] = 19.8873091
$row164[0,7] = 16.1697976
$row164[0,8] = 9.966300499999999
$row164[0,9] = 13.396861
$row164[0,10] = 16.8050699
$row164[0,11] = 28.7717943
$row164[0,12] = 30.6282619
$row164[0,14] = 11.1602497
$row164[0,15] = 27.3507073
$row164[0,16] = 29.102253
$row164[0,17] = 16.2861879
$row164[0,18] = 22.1573037
$row164[0,19] = 24.7486555
$row164[0,20] = 20.0155098
$row164[0,21] = 34.0333781
$row164[0,22] = 11.041463
$row164[0,23] = 14.2828296
$row164[0,24] = 10.9900457
$row164[0,25] = 15.0396092
$row164[0,26] = 18.4962669
$row164[0,27] = 21.8253941
$row164[0,29] = 33.1677202
$row164[0,30] = 19.2988576
$row164[0,31] = 20.6208636
$row164[0,32] = 21.0689236
$row164[0,33] = 21.7470439
$row164[0,34] = 10.2026457
$row164[0,35] = 11.8747435
$row164[0,36] = 17.6386579
$row164[0,37] = 23.367593
$row164[0,38] = 11.9330541
$row164[0,39] = 18.3314672
$row164[0,40] = 27.9180883
$row164[0,41] = 15.0514105
$row164[0,42] = 13.736239
$row164[0,44] = 11.389702
$row164[0,45] = 32.2731426
$row164[0,46] = 23.2163726
$row164[0,47] = 25.7752316
$row164[0,48] = 36.8295252
$row164[0,49] = 26.2702635
$row164[0,50] = 16.906638
$row164[0,52] = 8.4644455
$row164[0,53] = 16.9925478
$row164[0,54] = 19.5397754
$row164[0,55] = 19.4989579
$row164[0,56] = 19.7099731
$ws.Range("A164:BE164").Value = $row164

# Row 165
$row165 = New-Object 'object[,]' 1,57
$row165[0,0] = "13 07 2020"
$row165[0,1] = 21.0707577
$row165[0,2] = 35.0592101
$row165[0,3] = 25.6044353
$row165[0,5] = 33.9178388
$row165[0,6] = 20.1835229
$row165[0,7] = 16.2827761
$row165[0,8] = 10.2342095
$row165[0,9] = 12.962963
$row165[0,10] = 16.7743325
$row165[0,11] = 29.3247194
$row165[0,12] = 31.2932619
$row165[0,14] = 11.316568
$row165[0,15] = 27.87526
$row165[0,16] = 30.4004959
$row165[0,17] = 16.674027
$row165[0,18] = 22.1782766
$row165[0,19] = 25.2930529
$row165[0,20] = 20.4806406
$row165[0,21] = 34.3886942
$row165[0,22] = 11.0236526
$row165[0,23] = 14.2162917
$row165[0,24] = 10.8833242
$row165[0,25] = 15.0813544
$row165[0,26] = 19.2787578
$row165[0,27] = 22.2925693
$row165[0,29] = 33.0277497
$row165[0,30] = 20.0907005
$row165[0,31] = 20.6148735
$row165[0,32] = 21.0551163
$row165[0,33] = 22.4085361
$row165[0,34] = 10.3627672
$row165[0,35] = 11.8702682
$row165[0,36] = 17.8907546
$row165[0,37] = 24.1088139
$row165[0,38] = 11.8892303
$row165[0,39] = 18.677243
$row165[0,40] = 27.9416774
$row165[0,41] = 15.2937178
$row165[0,42] = 13.8130742
$row165[0,44] = 11.4429624
$row165[0,45] = 32.4104225
$row165[0,46] = 22.8237693
$row165[0,47] = 26.2592025
$row165[0,48] = 36.5479858
$row165[0,49] = 27.3497737
$row165[0,50] = 17.3724394
$row165[0,52] = 8.964498799999999
$row165[0,53] = 17.1777709
$row165[0,54] = 20.1971216
$row165[0,55] = 19.1533166
$row165[0,56] = 19.6826186
$ws.Range("A165:BE165").Value = $row165

# Row 166
$row166 = New-Object 'object[,]' 1,57
$row166[0,0] = "14 07 2020"
$row166[0,1] = 21.1761121
$row166[0,2] = 35.4318376
$row166[0,3] = 25.7271968
$row166[0,5] = 33.7522053
$row166[0,6] = 20.4414942
$row166[0,7] = 16.6002548
$row166[0,8] = 10.3587963
$row166[0,9] = 13.368669
$row166[0,10] = 17.0693028
$row166[0,11] = 29.6575934
$row166[0,12] = 31.6967999
$row166[0,14] = 11.3525391
$row166[0,15] = 27.8229325
$row166[0,16] = 31.2854616
$row166[0,17] = 16.8814154
$row166[0,18] = 22.7565397
$row166[0,19] = 25.4999925
$row166[0,20] = 20.702807
$row166[0,21] = 34.0504716
$row166[0,22] = 11.1837746
$row166[0,23] = 14.251715
$row166[0,24] = 10.543989
$row166[0,25] = 15.2867517
$row166[0,26] = 19.2937183
$row166[0,27] = 22.1773134
$row166[0,29] = 33.6324745
$row166[0,30] = 19.8004753
$row166[0,31] = 21.0336649
$row166[0,32] = 20.937093
$row166[0,33] = 22.4269818
$row166[0,34] = 10.1571899
$row166[0,35] = 12.2255893
$row166[0,36] = 18.6536533
$row166[0,37] = 24.0143003
$row166[0,38] = 11.6886585
$row166[0,39] = 19.1965769
$row166[0,40] = 28.1369815
$row166[0,41] = 15.0680596
$row166[0,42] = 14.0810271
$row166[0,44] = 11.2094101
$row166[0,45] = 32.7838938
$row166[0,46] = 24.5317403
$row166[0,47] = 26.5720371
$row166[0,48] = 36.6617892
$row166[0,49] = 26.8216569
$row166[0,50] = 17.6626752
$row166[0,52] = 8.6423404
$row166[0,53] = 17.3735058
$row166[0,54] = 20.3850363
$row166[0,55] = 19.1470154
$row166[0,56] = 20.4452796
$ws.Range("A166:BE166").Value = $row166

# Row 167
$row167 = New-Object 'object[,]' 1,57
$row167[0,0] = "15 07 2020"
$ws.Range("A167:BE167").Value = $row167
